$d = $word.ActiveDocument

# 1. Title: "Currículum vítae - Alex Wilber" -> "Currículum vítae: Alex Wilber"
$d.Content.Find.Execute("Currículum vítae - Alex Wilber", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Currículum vítae: Alex Wilber", 2)

# 2. Job title/dates: "Animación de Spark: Diseñador de animaciones (enero de 2021 - Presente)"
#    -> "Spark Animation: Diseñador de animaciones (enero de 2021 - actualidad)"
$d.Content.Find.Execute("Animación de Spark: Diseñador de animaciones (enero de 2021 - Presente)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Spark Animation: Diseñador de animaciones (enero de 2021 - actualidad)", 2)

# 3. Bullet: "Líder de un equipo de 12 animadores..." -> "Dirige un equipo de 12 animadores..."
#    (note: there is a non-breaking space between "12" and "animadores" in the source)
$d.Content.Find.Execute("Líder de un equipo de 12 animadores para crear animaciones 3D de alta calidad para diversos proyectos, como largometrajes, anuncios y videojuegos.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Dirige un equipo de 12 animadores para crear animaciones 3D de alta calidad para diversos proyectos, como largometrajes, anuncios y videojuegos.", 2)

# 4. "Pixel Studio: Diseñador de animaciones (junio de 2018 - dic 2020)"
#    -> "Pixel Studio: Diseñador de animaciones (junio de 2018 - diciembre de 2020)"
$d.Content.Find.Execute("Pixel Studio: Diseñador de animaciones (junio de 2018 - dic 2020)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Pixel Studio: Diseñador de animaciones (junio de 2018 - diciembre de 2020)", 2)

# 5. "Animación flash: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)"
#    -> "Flash Animation: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)"
$d.Content.Find.Execute("Animación flash: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Flash Animation: Diseñador de animaciones junior (septiembre de 2016 - mayo de 2018)", 2)

# 6. "Education" -> "Educación"
$d.Content.Find.Execute("Education", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Educación", 2)

# 7. "Maestro de Artes en Animación, Graduación esperada: dic 2025"
#    -> "Maestría en Humanidades con especialización en animación, graduación esperada: diciembre de 2025"
$d.Content.Find.Execute("Maestro de Artes en Animación, Graduación esperada: dic 2025", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Maestría en Humanidades con especialización en animación, graduación esperada: diciembre de 2025", 2)

# 8. "El arte de la animación 3D: una guía para principiantes." -> "The Art of 3D Animation: A Guide for Beginners."
$d.Content.Find.Execute("El arte de la animación 3D: una guía para principiantes.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The Art of 3D Animation: A Guide for Beginners.", 2)

Write-Host "done"
